# Excess mortality analyses - Week 50
# Update weekly observed/expected mortality figures for several existing
# weeks, extend the table with weeks 50-52, and move the totals row down
# to make room (old row 43 -> new row 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Revised observed mortality ("G" column) for weeks already present.
#    The "I" (Oversterfte = G-H) shared formula recalculates automatically.
# ---------------------------------------------------------------------
$ws.Range("G22").Value = 2673
$ws.Range("G23").Value = 2668
$ws.Range("G30").Value = 2719
$ws.Range("G32").Value = 2997
$ws.Range("G34").Value = 3216
$ws.Range("G35").Value = 3445
$ws.Range("G36").Value = 3675
$ws.Range("G38").Value = 3560
$ws.Range("G39").Value = 3317
$ws.Range("G40").Value = 3388
$ws.Range("G41").Value = 3494

# ---------------------------------------------------------------------
# 2) Extend the shared formula in column I down to the new row 42
#    (it used to stop at I41).
# ---------------------------------------------------------------------
$ws.Range("I42").Formula = "=G42-H42"

# ---------------------------------------------------------------------
# 3) Move the "Som week 11 tot en met 19" totals row from row 43 down to
#    row 46, freeing rows 42-44 for the new weekly data.
# ---------------------------------------------------------------------
$totalLabel = $ws.Range("F43").Value
$ws.Range("F43:I43").Clear()

# New week rows.
$ws.Range("F42").Value = 50
$ws.Range("G42").Value = 3571
$ws.Range("H42").Value = 3100

$ws.Range("F43").Value = 51
$ws.Range("F44").Value = 52

# Totals row, now on row 46.
$ws.Range("F46").Value = $totalLabel
$ws.Range("G46").NumberFormat = "0"
$ws.Range("H46").NumberFormat = "0"
$ws.Range("I46").NumberFormat = "0"
$ws.Range("G46").Formula = "=SUM(G3:G28)"
$ws.Range("H46").Formula = "=SUM(H3:H28)"
$ws.Range("I46").Formula = "=SUM(I3:I34)"

# ---------------------------------------------------------------------
# 4) Refresh the view: scroll position and active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I43").Select()
